# LunchModel_Ratings_Table: refresh pairwise "Price" comparison inputs and
# relabel the results header on the rating_table sheet.

$wb = $excel.ActiveWorkbook

$ratingScales = $wb.Worksheets.Item("rating_scales")
$ratingTable  = $wb.Worksheets.Item("rating_table")

# --- Pairwise "Price" comparison matrix updates (rating_scales) ---------
# 6-to-8-dollars vs 8-to-10-dollars
$ratingScales.Range("C12").Value = 0.9807692307692307
# 6-to-8-dollars vs more-than-10-dollars
$ratingScales.Range("D12").Value = 0.9622641509433962
# 8-to-10-dollars vs more-than-10-dollars
$ratingScales.Range("D13").Value = 0.9811320754716981

# --- Results header relabel (rating_table) -------------------------------
$ratingTable.Range("A7").Value = "ESTIMATED TOTALS AND PRIORITIES"

# --- Re-apply the rating dropdown validations -----------------------------
# Mirrors the source workbook re-registering its list validations (backed by
# the rating_scales lookup tables) for the rating input cells.
$validationTargets = @(
    @("B2:B4", "=rating_scales!A3:A6"),
    @("C2:C4", "=rating_scales!A12:A14"),
    @("D2:D4", "=rating_scales!A20:A23"),
    @("E2:E4", "=rating_scales!A29:A31")
)

foreach ($target in $validationTargets) {
    $addr = $target[0]
    $formula = $target[1]
    $ratingTable.Range($addr).Validation.Add(3, 1, 1, $formula)
}
